$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 1
$row2[0,1] = 0.3333333333333333
$row2[0,2] = 0.052073
$row2[0,3] = 0.156219
$row2[0,4] = 0.03816600682360385
$row2[0,5] = 0.03816600682360385
$row2[0,6] = 3
$row2[0,7] = 1
$row2[0,8] = 0.6692483333333333
$row2[0,9] = 2.007745
$row2[0,10] = 0.162783276239793
$row2[0,11] = 0.162783276239793
$row2[0,12] = 0.03484976846166667
$row2[0,13] = 0.313647916155
$row2[0,14] = 0.006212787631736532
$row2[0,15] = 0.006212787631736532
$ws.Range("E2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 1
$row3[0,1] = 0.3333333333333333
$row3[0,2] = 0.052073
$row3[0,3] = 0.156219
$row3[0,4] = 0.03816600682360385
$row3[0,5] = 0.03816600682360385
$row3[0,6] = 3
$row3[0,7] = 1
$row3[0,8] = 2.534005
$row3[0,9] = 7.602015
$row3[0,10] = 0.6163536244513372
$row3[0,11] = 0.6163536244513373
$row3[0,12] = 0.131953242365
$row3[0,13] = 1.187579181285
$row3[0,14] = 0.02352375663656271
$row3[0,15] = 0.02352375663656271
$ws.Range("E3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 1
$row4[0,1] = 0.3333333333333333
$row4[0,2] = 0.052073
$row4[0,3] = 0.156219
$row4[0,4] = 0.03816600682360385
$row4[0,5] = 0.03816600682360385
$row4[0,6] = 3
$row4[0,7] = 1
$row4[0,8] = 0.885104
$row4[0,9] = 2.655312
$row4[0,10] = 0.2152864964419472
$row4[0,11] = 0.2152864964419472
$row4[0,12] = 0.046090020592
$row4[0,13] = 0.414810185328
$row4[0,14] = 0.008216625892233124
$row4[0,15] = 0.008216625892233124
$ws.Range("E4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 1
$row5[0,1] = 0.3333333333333333
$row5[0,2] = 0.052073
$row5[0,3] = 0.156219
$row5[0,4] = 0.03816600682360385
$row5[0,5] = 0.03816600682360385
$row5[0,6] = 1
$row5[0,7] = 0.3333333333333333
$row5[0,8] = 0.022927
$row5[0,9] = 0.068781
$row5[0,10] = 0.005576602866922444
$row5[0,11] = 0.005576602866922445
$row5[0,12] = 0.001193877671
$row5[0,13] = 0.010744899039
$row5[0,14] = 0.0002128366630714908
$row5[0,15] = 0.0002128366630714909
$ws.Range("E5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3
$row6[0,1] = 1
$row6[0,2] = 0.7091626666666667
$row6[0,3] = 2.127488
$row6[0,4] = 0.5197685398391702
$row6[0,5] = 0.5197685398391702
$row6[0,6] = 3
$row6[0,7] = 1
$row6[0,8] = 0.6692483333333333
$row6[0,9] = 2.007745
$row6[0,10] = 0.162783276239793
$row6[0,11] = 0.162783276239793
$row6[0,12] = 0.4746059327288889
$row6[0,13] = 4.27145339456
$row6[0,14] = 0.0846096258013935
$row6[0,15] = 0.0846096258013935
$ws.Range("E6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3
$row7[0,1] = 1
$row7[0,2] = 0.7091626666666667
$row7[0,3] = 2.127488
$row7[0,4] = 0.5197685398391702
$row7[0,5] = 0.5197685398391702
$row7[0,6] = 3
$row7[0,7] = 1
$row7[0,8] = 2.534005
$row7[0,9] = 7.602015
$row7[0,10] = 0.6163536244513372
$row7[0,11] = 0.6163536244513373
$row7[0,12] = 1.797021743146667
$row7[0,13] = 16.17319568832
$row7[0,14] = 0.3203612234056518
$row7[0,15] = 0.3203612234056519
$ws.Range("E7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 3
$row8[0,1] = 1
$row8[0,2] = 0.7091626666666667
$row8[0,3] = 2.127488
$row8[0,4] = 0.5197685398391702
$row8[0,5] = 0.5197685398391702
$row8[0,6] = 3
$row8[0,7] = 1
$row8[0,8] = 0.885104
$row8[0,9] = 2.655312
$row8[0,10] = 0.2152864964419472
$row8[0,11] = 0.2152864964419472
$row8[0,12] = 0.6276827129173334
$row8[0,13] = 5.649144416256
$row8[0,14] = 0.1118991479027216
$row8[0,15] = 0.1118991479027216
$ws.Range("E8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 3
$row9[0,1] = 1
$row9[0,2] = 0.7091626666666667
$row9[0,3] = 2.127488
$row9[0,4] = 0.5197685398391702
$row9[0,5] = 0.5197685398391702
$row9[0,6] = 1
$row9[0,7] = 0.3333333333333333
$row9[0,8] = 0.022927
$row9[0,9] = 0.068781
$row9[0,10] = 0.005576602866922444
$row9[0,11] = 0.005576602866922445
$row9[0,12] = 0.01625897245866667
$row9[0,13] = 0.146330752128
$row9[0,14] = 0.002898542729403209
$row9[0,15] = 0.002898542729403209
$ws.Range("E9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 2
$row10[0,1] = 0.6666666666666666
$row10[0,2] = 0.299804
$row10[0,3] = 0.8994119999999999
$row10[0,4] = 0.2197361686429383
$row10[0,5] = 0.2197361686429384
$row10[0,6] = 3
$row10[0,7] = 1
$row10[0,8] = 0.6692483333333333
$row10[0,9] = 2.007745
$row10[0,10] = 0.162783276239793
$row10[0,11] = 0.162783276239793
$row10[0,12] = 0.2006433273266666
$row10[0,13] = 1.80578994594
$row10[0,14] = 0.03576937344007718
$row10[0,15] = 0.03576937344007718
$ws.Range("E10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,16
$row11[0,0] = 2
$row11[0,1] = 0.6666666666666666
$row11[0,2] = 0.299804
$row11[0,3] = 0.8994119999999999
$row11[0,4] = 0.2197361686429383
$row11[0,5] = 0.2197361686429384
$row11[0,6] = 3
$row11[0,7] = 1
$row11[0,8] = 2.534005
$row11[0,9] = 7.602015
$row11[0,10] = 0.6163536244513372
$row11[0,11] = 0.6163536244513373
$row11[0,12] = 0.75970483502
$row11[0,13] = 6.837343515179999
$row11[0,14] = 0.1354351839661253
$row11[0,15] = 0.1354351839661254
$ws.Range("E11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,16
$row12[0,0] = 2
$row12[0,1] = 0.6666666666666666
$row12[0,2] = 0.299804
$row12[0,3] = 0.8994119999999999
$row12[0,4] = 0.2197361686429383
$row12[0,5] = 0.2197361686429384
$row12[0,6] = 3
$row12[0,7] = 1
$row12[0,8] = 0.885104
$row12[0,9] = 2.655312
$row12[0,10] = 0.2152864964419472
$row12[0,11] = 0.2152864964419472
$row12[0,12] = 0.265357719616
$row12[0,13] = 2.388219476544
$row12[0,14] = 0.04730622988871505
$row12[0,15] = 0.04730622988871507
$ws.Range("E12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,16
$row13[0,0] = 2
$row13[0,1] = 0.6666666666666666
$row13[0,2] = 0.299804
$row13[0,3] = 0.8994119999999999
$row13[0,4] = 0.2197361686429383
$row13[0,5] = 0.2197361686429384
$row13[0,6] = 1
$row13[0,7] = 0.3333333333333333
$row13[0,8] = 0.022927
$row13[0,9] = 0.068781
$row13[0,10] = 0.005576602866922444
$row13[0,11] = 0.005576602866922445
$row13[0,12] = 0.006873606307999999
$row13[0,13] = 0.06186245677199999
$row13[0,14] = 0.001225381348020764
$row13[0,15] = 0.001225381348020764
$ws.Range("E13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,16
$row14[0,0] = 2
$row14[0,1] = 0.6666666666666666
$row14[0,2] = 0.303342
$row14[0,3] = 0.910026
$row14[0,4] = 0.2223292846942876
$row14[0,5] = 0.2223292846942876
$row14[0,6] = 3
$row14[0,7] = 1
$row14[0,8] = 0.6692483333333333
$row14[0,9] = 2.007745
$row14[0,10] = 0.162783276239793
$row14[0,11] = 0.162783276239793
$row14[0,12] = 0.20301112793
$row14[0,13] = 1.82710015137
$row14[0,14] = 0.03619148936658581
$row14[0,15] = 0.03619148936658581
$ws.Range("E14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,16
$row15[0,0] = 2
$row15[0,1] = 0.6666666666666666
$row15[0,2] = 0.303342
$row15[0,3] = 0.910026
$row15[0,4] = 0.2223292846942876
$row15[0,5] = 0.2223292846942876
$row15[0,6] = 3
$row15[0,7] = 1
$row15[0,8] = 2.534005
$row15[0,9] = 7.602015
$row15[0,10] = 0.6163536244513372
$row15[0,11] = 0.6163536244513373
$row15[0,12] = 0.7686701447100001
$row15[0,13] = 6.91803130239
$row15[0,14] = 0.1370334604429974
$row15[0,15] = 0.1370334604429974
$ws.Range("E15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,16
$row16[0,0] = 2
$row16[0,1] = 0.6666666666666666
$row16[0,2] = 0.303342
$row16[0,3] = 0.910026
$row16[0,4] = 0.2223292846942876
$row16[0,5] = 0.2223292846942876
$row16[0,6] = 3
$row16[0,7] = 1
$row16[0,8] = 0.885104
$row16[0,9] = 2.655312
$row16[0,10] = 0.2152864964419472
$row16[0,11] = 0.2152864964419472
$row16[0,12] = 0.268489217568
$row16[0,13] = 2.416402958112
$row16[0,14] = 0.04786449275827741
$row16[0,15] = 0.04786449275827743
$ws.Range("E16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,16
$row17[0,0] = 2
$row17[0,1] = 0.6666666666666666
$row17[0,2] = 0.303342
$row17[0,3] = 0.910026
$row17[0,4] = 0.2223292846942876
$row17[0,5] = 0.2223292846942876
$row17[0,6] = 1
$row17[0,7] = 0.3333333333333333
$row17[0,8] = 0.022927
$row17[0,9] = 0.068781
$row17[0,10] = 0.005576602866922444
$row17[0,11] = 0.005576602866922445
$row17[0,12] = 0.006954722034
$row17[0,13] = 0.062592498306
$row17[0,14] = 0.001239842126426981
$row17[0,15] = 0.001239842126426981
$ws.Range("E17:T17").Value = $row17
